$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (shifts old rows 85:136 down to 86:137)
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with the inserted record's data
$ws.Range("A85").Value = 4
$ws.Range("B85").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C85").Value = "Los Lagos"
$ws.Range("D85").Value = 45097
$ws.Range("E85").Value = 10
$ws.Range("F85").Value = 100112026
$ws.Range("G85").Value = "Haba"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 80
$ws.Range("K85").Value = 21000
$ws.Range("L85").Value = 21000
$ws.Range("M85").Value = 21000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Provincia de Limarí"
$ws.Range("P85").Value = 840
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"

Write-Host "done"
